# Add support for a reed switch: new register row (row 10) describing the
# "reed" field (R[7:0]=0x00 magnetic/L, R=0xAA released/H), one byte wide.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 10 data: size=1 byte, type="u", field="reed", description, and
# one "R" marker per bit-column (F..U) to show which byte it occupies.
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "u"
$ws.Range("D10").Value = "reed"
$ws.Range("E10").Value = "reed switch state: R[7:0]=0x00 (magnetic, L), R=0xAA (released, H)"
$ws.Range("F10:U10").Value = "R"

# The edit's last user action left the selection on E11.
$ws.Range("E11").Select() | Out-Null
